$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit reshuffles the per-row market data (date, variety, quality,
# volume, prices, unit, origin, $/kg, kg-or-units) across rows 2-30 --
# each destination row ends up with the values that used to live on a
# different source row. Snapshot every source cell via COM first (one
# read per cell) so the subsequent writes do not clobber data we still
# need to copy elsewhere.
$cols = @("D","H","I","J","K","L","M","N","O","P","Q")

$mapping = @{}
$mapping[2] = 7
$mapping[3] = 23
$mapping[4] = 24
$mapping[5] = 25
$mapping[6] = 26
$mapping[7] = 5
$mapping[8] = 12
$mapping[9] = 8
$mapping[10] = 10
$mapping[11] = 11
$mapping[12] = 22
$mapping[13] = 15
$mapping[14] = 18
$mapping[15] = 3
$mapping[16] = 2
$mapping[17] = 13
$mapping[18] = 14
$mapping[19] = 16
$mapping[20] = 9
$mapping[21] = 6
$mapping[22] = 19
$mapping[23] = 30
$mapping[24] = 20
$mapping[25] = 21
$mapping[26] = 4
$mapping[27] = 27
$mapping[28] = 28
$mapping[29] = 29
$mapping[30] = 17

# Only source rows actually referenced by the mapping need to be read.
$sourceRows = $mapping.Values | Sort-Object -Unique

$snapshot = @{}
foreach ($r in $sourceRows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Now write each destination row using the snapshot taken above.
foreach ($r in ($mapping.Keys | Sort-Object)) {
    $src = $snapshot[$mapping[$r]]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $src[$col]
    }
}
